$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 10 new rows starting at row 142, pushing existing row 142 (train_test_80_20) and below down by 10 rows.
$insertRange = $ws.Range("A142:A151")
$insertRange.EntireRow.Insert()

# New values to place in rows 142-152 (replacing old row142 "train_test_80_20" with 11 new rows)
$newValues = @(
    "Altman_Z_diff",
    "Ratio_A_diff",
    "Ratio_B_diff",
    "Ratio_C_diff",
    "Ratio_D_diff",
    "Ratio_E_diff",
    "grossProfitRatio_diff",
    "ebitdaratio_diff",
    "operatingIncomeRatio_diff",
    "incomeBeforeTaxRatio_diff",
    "netIncomeRatio_diff"
)

$row = 142
foreach ($val in $newValues) {
    $ws.Cells.Item($row, 1).Value = $val
    $row = $row + 1
}

# Append two new rows at the end (now rows 171 and 172)
$ws.Cells.Item(171, 1).Value = "num_q_by_len"
$ws.Cells.Item(172, 1).Value = "pos_score_finbert"
